# Edit script: "Documentation updated with redux with angular"
$d = $word.ActiveDocument

# --- Change 1: "Builds" heading paragraph ---
# Drop <w:lang w:val="en-US"/> from the paragraph mark and run
# properties; wrap the run text in spell-check proofErr markers.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML('<w:p w:rsidR="00B74379" w:rsidRDefault="00B74379"><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Builds</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# --- Change 2: "--allow-unrelated-histories" paragraph ---
# Remove the _GoBack bookmark here; it is relocated to the very end
# of the document (see change 4 below).
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML('<w:p w:rsidR="00B74379" w:rsidRDefault="00B74379"><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Para ello necesitamos añadir el parámetro </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="242729"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="EFF0F1"/></w:rPr><w:t>--allow-unrelated-histories</w:t></w:r></w:p>')

# --- Change 3: "ASP.NET Core Auth" heading paragraph ---
# Add <w:lang w:val="en-US"/> to the paragraph mark and both runs.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML('<w:p w:rsidR="004A3615" w:rsidRPr="00B74379" w:rsidRDefault="004A3615"><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00B74379"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">ASP.NET Core </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B74379"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Auth</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# --- Change 4: new "Redux with Angular" section ---
# The last empty (sz 20) paragraph before the document's final
# paragraph is replaced with the new heading, description, npm
# command and link paragraphs; the _GoBack bookmark now sits at the
# very end of this new content.
$target = $d.Paragraphs.Item(59)
$target.Range.InsertXML('<w:p w:rsidR="00E735FB" w:rsidRPr="00E735FB" w:rsidRDefault="00E735FB" w:rsidP="00E735FB"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Redux with Angular</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Esta es la librería para usar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>redux</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> con angular. (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>incluye</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> como dependencia </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Redux</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:spacing w:val="-1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:spacing w:val="-1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">$ </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:spacing w:val="-1"/><w:lang w:val="en-US"/></w:rPr><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:spacing w:val="-1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> install redux @angular-redux/store —save</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:spacing w:val="-1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>https://medium.com/codingthesmartway-com-blog/angular-and-redux-ecd22ea53492</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
